# Updated symbol list on Wed Dec 14 23:49:06 UTC 2022 with GitHub Actions
# Refresh the "Price" column (D) values for the crypto rows whose quotes moved.
# Values are stored as text in the sheet (t="inlineStr"/shared string), so we
# write them with a leading apostrophe to force text entry and then reset the
# cell style back to "Normal" so no stray number-format/style is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'268.44"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'22.83"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'6.304"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'0.06187"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'3.581"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'6.700"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'1.366"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.8397"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.01363"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.1604"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.08246"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'0.03259"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'3.901"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'0.001714"
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'0.006269"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.005362"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'3.762"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'0.1212"
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'0.0002683"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.04676"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'0.006964"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'0.1151"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.003601"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.01208"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.00006228"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'0.00000000750"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'0.7002"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'0.1652"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'0.00002101"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'0.01240"
$c.Style = "Normal"
